# Generate Report for Handoff
#
# The CI run that produced this localization-status report re-ran the
# handback/handoff pass later than the previous commit; several files that
# are still "Ready for handoff" (or hit "Handback transform failed") now
# share one single, later timestamp (the most recent run time) instead of
# each carrying its own older timestamp. This updates the Overview sheet's
# "Latest Handoff Date" column and each locale sheet's "Latest Handoff
# Datetime" column for those rows.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D ("Latest Handoff Date") ---
$ws = $wb.Worksheets.Item("Overview")
foreach ($r in 7,9,10,11,12,13,14,15,16) {
    $ws.Cells.Item($r, 4).Value = "2016-26-12 04:26:42"
}

# --- zh-cn sheet: column E ("Latest Handoff Datetime") ---
$ws = $wb.Worksheets.Item("zh-cn")
foreach ($r in 7,9,10,11,12,13,14,15,16) {
    $ws.Cells.Item($r, 5).Value = "2016-03-12 04:26:38"
}

# --- de-de sheet: column E ("Latest Handoff Datetime") ---
$ws = $wb.Worksheets.Item("de-de")
foreach ($r in 7,9,10,11,12,13,14,15,16) {
    $ws.Cells.Item($r, 5).Value = "2016-03-12 04:26:42"
}
